# Marks A' level sheet SUBMATH - add "Mid Paper 1" column (E) and enter
# mid-paper percentage marks for the two students who already have a
# Paper 1 score, on both the "Senior Six" and "Senior Five" sheets, then
# make "Senior Six" the active/selected sheet (it was "Senior Five").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Senior Six"
$ws2 = $wb.Worksheets.Item(2)   # "Senior Five"

foreach ($ws in @($ws1, $ws2)) {
    # New column header
    $ws.Range("E1").Value = "Mid Paper 1"

    # Give the previously-empty E3/E4 cells an explicit (blank) style so
    # they line up with the rest of the formatted row, matching C3/C4.
    $ws.Cells.Item(3, 5).Style = "Normal"
    $ws.Cells.Item(4, 5).Style = "Normal"

    # Mid Paper 1 marks for the two graded students
    $ws.Range("E5").Value = 70
    $ws.Range("E6").Value = 71

    # Move the selection over to the new column
    $ws.Range("E1:E6").Select()
}

# "Senior Six" is now the active/selected sheet (was "Senior Five")
$ws1.Activate()
